$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: quarter period headers (shift left, append new Q4 1401/12) ---
$ws.Range("D8").Value = 'فصل سوم منتهی به 1399/09'
$ws.Range("E8").Value = 'فصل چهارم منتهی به 1399/12'
$ws.Range("F8").Value = 'فصل اول منتهی به 1400/03'
$ws.Range("G8").Value = 'فصل دوم منتهی به 1400/06'
$ws.Range("H8").Value = 'فصل سوم منتهی به 1400/09'
$ws.Range("I8").Value = 'فصل چهارم منتهی به 1400/12'
$ws.Range("J8").Value = 'فصل اول منتهی به 1401/03'
$ws.Range("K8").Value = 'فصل دوم منتهی به 1401/06'
$ws.Range("L8").Value = 'فصل سوم منتهی به 1401/09'
$ws.Range("M8").Value = 'فصل چهارم منتهی به 1401/12'

# --- Row 9: publish date headers (shift left, append new date) ---
# Note: a couple of these look like real calendar dates (YYYY-MM-DD) and Excel
# auto-converts them to date serials; force them to text and restore formatting.
$ws.Range("D9").Value = '1400-10-29 (3)'
$ws.Range("E9").Value = '1401-03-11 (8)'
$ws.Range("F9").Value = '1401-04-29 (2)'
$ws.Range("G9").Value = '1401-08-29 (4)'
$ws.Range("H9").Value = '1401-10-28 (2)'
$ws.Range("I9").Value = '1402-02-30 (7)'
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = '1401-04-29'
$ws.Range("K9").Value = '1401-08-29 (2)'
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = '1401-10-28'
$ws.Range("M9").Value = '1402-02-30'
# restore original cell formatting/style (border etc.) on the two forced-text cells
$ws.Range("D9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Data rows 12-56: each quarter column shifts left by one, newest quarter (M) filled in ---
# row 12
$ws.Range("D12").Value = 215556
$ws.Range("E12").Value = -129408
$ws.Range("F12").Value = -632852
$ws.Range("G12").Value = 1474211
$ws.Range("H12").Value = -903530
$ws.Range("I12").Value = 1042157
$ws.Range("J12").Value = 91306
$ws.Range("K12").Value = 789983
$ws.Range("L12").Value = 173689
$ws.Range("M12").Value = 2075749

# row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = -5083
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = -2456
$ws.Range("K13").Value = -42430
$ws.Range("L13").Value = -48000
$ws.Range("M13").Value = -98535

# row 14
$ws.Range("D14").Value = 215556
$ws.Range("E14").Value = -134491
$ws.Range("F14").Value = -632852
$ws.Range("G14").Value = 1474211
$ws.Range("H14").Value = -903530
$ws.Range("I14").Value = 1042157
$ws.Range("J14").Value = 88850
$ws.Range("K14").Value = 747553
$ws.Range("L14").Value = 125689
$ws.Range("M14").Value = 1977214

# row 16
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0

# row 17
$ws.Range("D17").Value = -47606
$ws.Range("E17").Value = -108366
$ws.Range("F17").Value = -69081
$ws.Range("G17").Value = -289180
$ws.Range("H17").Value = -187429
$ws.Range("I17").Value = -96272
$ws.Range("J17").Value = -75435
$ws.Range("K17").Value = -111903
$ws.Range("L17").Value = -136003
$ws.Range("M17").Value = -362289

# row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0

# row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0

# row 20
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0

# row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# row 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

# row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# row 24
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0

# row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 9989
$ws.Range("F25").Value = 9893
$ws.Range("G25").Value = 321
$ws.Range("H25").Value = -2
$ws.Range("I25").Value = -9979
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0

# row 26
$ws.Range("D26").Value = -52
$ws.Range("E26").Value = 55
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = -3
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = -7
$ws.Range("M26").Value = 0

# row 27
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0

# row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0

# row 29
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0

# row 30
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 592
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 387
$ws.Range("H30").Value = 31
$ws.Range("I30").Value = 185
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 136
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 289

# row 31
$ws.Range("D31").Value = 372
$ws.Range("E31").Value = 3470
$ws.Range("F31").Value = 72750
$ws.Range("G31").Value = -67059
$ws.Range("H31").Value = 126741
$ws.Range("I31").Value = -123913
$ws.Range("J31").Value = 1010
$ws.Range("K31").Value = 1063
$ws.Range("L31").Value = 358
$ws.Range("M31").Value = 321

# row 32
$ws.Range("D32").Value = -47286
$ws.Range("E32").Value = -94260
$ws.Range("F32").Value = 13562
$ws.Range("G32").Value = -355531
$ws.Range("H32").Value = -60659
$ws.Range("I32").Value = -229979
$ws.Range("J32").Value = -74428
$ws.Range("K32").Value = -110701
$ws.Range("L32").Value = -135652
$ws.Range("M32").Value = -361679

# row 33
$ws.Range("D33").Value = 168270
$ws.Range("E33").Value = -228751
$ws.Range("F33").Value = -619290
$ws.Range("G33").Value = 1118680
$ws.Range("H33").Value = -964189
$ws.Range("I33").Value = 812178
$ws.Range("J33").Value = 14422
$ws.Range("K33").Value = 636852
$ws.Range("L33").Value = -9963
$ws.Range("M33").Value = 1615535

# row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 187018
$ws.Range("F35").Value = 15658
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = -1
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0

# row 36
$ws.Range("D36").Value = '-'
$ws.Range("E36").Value = '-'
$ws.Range("F36").Value = '-'
$ws.Range("G36").Value = '-'
$ws.Range("H36").Value = '-'
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = '-'
$ws.Range("K36").Value = '-'
$ws.Range("L36").Value = '-'
$ws.Range("M36").Value = 0

# row 37
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0

# row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0

# row 39
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 550000
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 296000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0

# row 40
$ws.Range("D40").Value = -76847
$ws.Range("E40").Value = -52520
$ws.Range("F40").Value = -69622
$ws.Range("G40").Value = -134749
$ws.Range("H40").Value = 158439
$ws.Range("I40").Value = -490462
$ws.Range("J40").Value = -156579
$ws.Range("K40").Value = -74742
$ws.Range("L40").Value = -78285
$ws.Range("M40").Value = 0

# row 41
$ws.Range("D41").Value = -3118
$ws.Range("E41").Value = -634
$ws.Range("F41").Value = -10530
$ws.Range("G41").Value = -12414
$ws.Range("H41").Value = -9758
$ws.Range("I41").Value = -14886
$ws.Range("J41").Value = -8688
$ws.Range("K41").Value = -5555
$ws.Range("L41").Value = -2013
$ws.Range("M41").Value = 0

# row 42
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0

# row 43
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0

# row 44
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0

# row 45
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0

# row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0

# row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0

# row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0

# row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0

# row 50
$ws.Range("D50").Value = -5766
$ws.Range("E50").Value = -4888
$ws.Range("F50").Value = -110
$ws.Range("G50").Value = -144
$ws.Range("H50").Value = -32998
$ws.Range("I50").Value = -177374
$ws.Range("J50").Value = -20
$ws.Range("K50").Value = -20236
$ws.Range("L50").Value = -150116
$ws.Range("M50").Value = -758929

# row 51
$ws.Range("D51").Value = -85731
$ws.Range("E51").Value = 128976
$ws.Range("F51").Value = 485396
$ws.Range("G51").Value = -147307
$ws.Range("H51").Value = 115683
$ws.Range("I51").Value = -386723
$ws.Range("J51").Value = -165287
$ws.Range("K51").Value = -100533
$ws.Range("L51").Value = -230414
$ws.Range("M51").Value = -758929

# row 52
$ws.Range("D52").Value = 82539
$ws.Range("E52").Value = -99775
$ws.Range("F52").Value = -133894
$ws.Range("G52").Value = 971373
$ws.Range("H52").Value = -848506
$ws.Range("I52").Value = 425455
$ws.Range("J52").Value = -150865
$ws.Range("K52").Value = 536319
$ws.Range("L52").Value = -240377
$ws.Range("M52").Value = 856606

# row 53
$ws.Range("D53").Value = 405852
$ws.Range("E53").Value = 482234
$ws.Range("F53").Value = 379917
$ws.Range("G53").Value = 246023
$ws.Range("H53").Value = 1217219
$ws.Range("I53").Value = 368712
$ws.Range("J53").Value = 803430
$ws.Range("K53").Value = 652972
$ws.Range("L53").Value = 1189394
$ws.Range("M53").Value = 954083

# row 54
$ws.Range("D54").Value = -6157
$ws.Range("E54").Value = 7619
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = -177
$ws.Range("H54").Value = -1
$ws.Range("I54").Value = -898
$ws.Range("J54").Value = 407
$ws.Range("K54").Value = 103
$ws.Range("L54").Value = 5066
$ws.Range("M54").Value = 15259

# row 55
$ws.Range("D55").Value = 482234
$ws.Range("E55").Value = 390078
$ws.Range("F55").Value = 246023
$ws.Range("G55").Value = 1217219
$ws.Range("H55").Value = 368712
$ws.Range("I55").Value = 803430
$ws.Range("J55").Value = 652972
$ws.Range("K55").Value = 1189394
$ws.Range("L55").Value = 954083
$ws.Range("M55").Value = 1825948

# row 56
$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 147324
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0

# --- Column width adjustments: the "wide" (31) column marking fiscal year-end shifts ---
# from columns F & J to columns E, I and the new M (year-end Q4 1401/12)
$ws.Columns.Item(5).ColumnWidth = 30.2   # E -> 31
$ws.Columns.Item(6).ColumnWidth = 28.2   # F -> 29
$ws.Columns.Item(9).ColumnWidth = 30.2   # I -> 31
$ws.Columns.Item(10).ColumnWidth = 28.2  # J -> 29
$ws.Columns.Item(13).ColumnWidth = 30.2  # M -> 31
